# Refitting NCDEs to individual patients (for manuscript figure)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" header in column H, matching the header style used by the
# other header cells (bold, centered, bordered) via format-only copy/paste.
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Update refitted Prediction/Error (and Accuracy for the last patient) values,
# and populate the new Label column (0 = Control, 1 = MDD) for every patient row.

# --- Control 1 (batchsize 100) iteration block ---
$ws.Range("D2").Value = 0.5479512342500286
$ws.Range("E2").Value = 0.5479512342500286
$ws.Range("H2").Value = 0

$ws.Range("H3").Value = 0

$ws.Range("D4").Value = 0.4103447143359923
$ws.Range("E4").Value = 0.4103447143359923
$ws.Range("H4").Value = 0

$ws.Range("D5").Value = 0.3444163388797605
$ws.Range("E5").Value = 0.3444163388797605
$ws.Range("H5").Value = 0

$ws.Range("D6").Value = 0.4571451606035604
$ws.Range("E6").Value = 0.4571451606035604
$ws.Range("H6").Value = 0

$ws.Range("D7").Value = 0.6721420294334516
$ws.Range("E7").Value = 0.3278579705665484
$ws.Range("H7").Value = 1

$ws.Range("D8").Value = 0.5175788968357256
$ws.Range("E8").Value = 0.4824211031642744
$ws.Range("H8").Value = 1

$ws.Range("D9").Value = 0.4199118772871552
$ws.Range("E9").Value = 0.5800881227128448
$ws.Range("H9").Value = 1

$ws.Range("D10").Value = 0.4677700188085754
$ws.Range("E10").Value = 0.5322299811914246
$ws.Range("H10").Value = 1

$ws.Range("D11").Value = 0.4516828654224521
$ws.Range("E11").Value = 0.5483171345775479
$ws.Range("F11").Value = 0.6385025978088379
$ws.Range("H11").Value = 1

# --- Control 1 (batchsize 200) iteration block ---
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
